# Se agrego historia de usuario HUFP-005 en el documento Peaceathome.docx
# dentro de la rama fpaiz-2024264
#
# Adds a new "HUFP-005: Registro de nuevos empleados" user-story block
# right after the existing HUFP-004 block, reusing the (hidden) _GoBack
# bookmark the way the original author's edit did.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Locate the paragraph that ends the HUFP-004 block ("Para: Poder
#    gestionar de manera eficiente ...").  We walk the Paragraphs
#    collection explicitly instead of relying on Find, since we need the
#    paragraph INDEX for later Item() lookups (Paragraph.Next is not
#    reliable in this host).
# ---------------------------------------------------------------------
$total = $d.Paragraphs.Count
$paraIdx = -1
for ($i = 1; $i -le $total; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*Poder gestionar de manera eficiente*") {
        $paraIdx = $i
        break
    }
}

# The blank separator paragraph that already exists right after it.
$blankIdx = $paraIdx + 1
$blankPara = $d.Paragraphs.Item($blankIdx)

# ---------------------------------------------------------------------
# 2. Insert 6 new paragraphs right after that blank separator:
#      heading  (HUFP-005: Registro de nuevos empleados)
#      Como
#      Quiero
#      Para
#      blank
#      blank
#    (the document ends up with two trailing blank paragraphs before the
#    final Arial-styled paragraph, matching the target diff.)
#    Calling InsertParagraphAfter() repeatedly on the SAME blank
#    paragraph's Range always inserts right after it, pushing the
#    previously inserted ones further down - so six calls yield six
#    new paragraphs in the right order.
# ---------------------------------------------------------------------
for ($n = 0; $n -lt 6; $n++) {
    $blankPara.Range.InsertParagraphAfter()
}

$headingIdx = $blankIdx + 1
$comoIdx    = $blankIdx + 2
$quieroIdx  = $blankIdx + 3
$paraIdx2   = $blankIdx + 4
# ($blankIdx + 5) and ($blankIdx + 6) are the two new trailing blank
# paragraphs - left empty.

# Word color value for OOXML color 0070C0 (stored BGR: 0x00C07000).
$blue = 12611584

# =======================================================================
# 3. Heading paragraph: "HUFP -005: Registro de nuevos empleados"
#    "HUFP -005:" bold+blue+underline, rest bold+blue (no underline).
#    The bookmark _GoBack wraps "HUFP -005" (without the colon).
# =======================================================================
$headingPara = $d.Paragraphs.Item($headingIdx)
$headingText = "HUFP -005: Registro de nuevos empleados"
$headingPara.Range.InsertAfter($headingText)

$headingRange = $d.Paragraphs.Item($headingIdx).Range
$headingStart = $headingRange.Start

# Whole paragraph -> bold + blue.
$headingRange.Font.Bold = $true
$headingRange.Font.Color = $blue

# "HUFP -005:" (first 11 characters) -> also underlined.
$underlineLen = ("HUFP -005:").Length
$underlineRange = $d.Range($headingStart, $headingStart + $underlineLen)
$underlineRange.Font.Underline = 1

# Bookmark _GoBack wraps "HUFP -005" only (9 chars, excludes the colon).
# Re-using the name "_GoBack" relocates the document's existing hidden
# bookmark instead of creating a duplicate.
$bookmarkLen = ("HUFP -005").Length
$bookmarkRange = $d.Range($headingStart, $headingStart + $bookmarkLen)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# =======================================================================
# 4. "Como: Responsable de contratación"
# =======================================================================
$comoPara = $d.Paragraphs.Item($comoIdx)
$comoPara.Range.InsertAfter("Como: Responsable de contratación")
$comoRange = $d.Paragraphs.Item($comoIdx).Range
$comoLabelRange = $d.Range($comoRange.Start, $comoRange.Start + ("Como:").Length)
$comoLabelRange.Font.Bold = $true

# =======================================================================
# 5. "Quiero: Poder agregar nuevos empleados con su información
#     (nombre, apellido, correo, teléfono, cargo, fecha de contratación)"
# =======================================================================
$quieroPara = $d.Paragraphs.Item($quieroIdx)
$quieroText = "Quiero: Poder agregar nuevos empleados con su información (nombre, apellido, correo, teléfono, cargo, fecha de contratación)"
$quieroPara.Range.InsertAfter($quieroText)
$quieroRange = $d.Paragraphs.Item($quieroIdx).Range
$quieroLabelRange = $d.Range($quieroRange.Start, $quieroRange.Start + ("Quiero:").Length)
$quieroLabelRange.Font.Bold = $true

# =======================================================================
# 6. "Para: Mantener la base de datos actualizada y asegurarme de que
#     se realice un registro adecuado y completo de todos los empleados
#     nuevos."
# =======================================================================
$paraPara2 = $d.Paragraphs.Item($paraIdx2)
$paraText2 = "Para: Mantener la base de datos actualizada y asegurarme de que se realice un registro adecuado y completo de todos los empleados nuevos."
$paraPara2.Range.InsertAfter($paraText2)
$paraRange2 = $d.Paragraphs.Item($paraIdx2).Range
$paraLabelRange = $d.Range($paraRange2.Start, $paraRange2.Start + ("Para:").Length)
$paraLabelRange.Font.Bold = $true

Write-Host "Done. Total paragraphs: $($d.Paragraphs.Count)"
